$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 5045357.142857142
$ws.Range("C3").Value = 3559285.714285714
$ws.Range("C4").Value = 2829642.857142857
$ws.Range("C5").Value = 6202857.142857143
$ws.Range("C6").Value = 2900714.285714285
$ws.Range("C7").Value = 20537857.14285714
